$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '72.695.74'
$ws.Range("E2").Value = '  +5.63%  '
$ws.Range("D3").Value = '4.061.19'
$ws.Range("E3").Value = '  +5.52%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '''522.46'
$ws.Range("E5").Value = '  -0.19%  '
$ws.Range("D6").Value = '''148.24'
$ws.Range("E6").Value = '  +4.24%  '
$ws.Range("D7").Value = '''0.730'
$ws.Range("E7").Value = '  +20.41%  '
$ws.Range("D8").Value = '4.051.96'
$ws.Range("E8").Value = '  +5.52%  '
$ws.Range("E9").Value = '  +0.14%  '
$ws.Range("E10").Value = '  +10.41%  '
$ws.Range("E11").Value = '  +4.75%  '
$ws.Range("D12").Value = '''0.0000332'
$ws.Range("E12").Value = '  +0.93%  '
$ws.Range("D13").Value = '''48.50'
$ws.Range("E13").Value = '  +16.59%  '
$ws.Range("D14").Value = '''11.14'
$ws.Range("E14").Value = '  +10.16%  '
$ws.Range("D15").Value = '4.703.96'
$ws.Range("E15").Value = '  +5.35%  '
$ws.Range("D16").Value = '4.061.33'
$ws.Range("E16").Value = '  +4.40%  '
$ws.Range("E17").Value = '  +4.39%  '
$ws.Range("D18").Value = '''14.40'
$ws.Range("E18").Value = '  +3.95%  '
$ws.Range("E19").Value = '  +1.88%  '
$ws.Range("E20").Value = '  -0.14%  '
$ws.Range("D21").Value = '72.597.20'
$ws.Range("E21").Value = '  +5.55%  '
$ws.Range("D22").Value = '''453.59'
$ws.Range("E22").Value = '  +8.00%  '
$ws.Range("D23").Value = '''105.42'
$ws.Range("E23").Value = '  +21.39%  '
$ws.Range("D24").Value = '''3.61'
$ws.Range("E24").Value = '  +7.09%  '
$ws.Range("E25").Value = '  +8.02%  '
$ws.Range("E26").Value = '  +2.20%  '
$ws.Range("D27").Value = '''11.40'
$ws.Range("E27").Value = '  +0.97%  '
$ws.Range("D28").Value = '''11.10'
$ws.Range("E28").Value = '  +5.79%  '
$ws.Range("D29").Value = '''38.23'
$ws.Range("E29").Value = '  +6.43%  '
$ws.Range("E30").Value = '  +3.02%  '
$ws.Range("E31").Value = '  +16.79%  '
$ws.Range("E32").Value = '  +4.96%  '
$ws.Range("E33").Value = '  +4.58%  '
$ws.Range("E34").Value = '  -0.90%  '
$ws.Range("D35").Value = '''67.78'
$ws.Range("E35").Value = '  -0.30%  '
$ws.Range("D36").Value = '''6.62'
$ws.Range("E36").Value = '  +12.56%  '
$ws.Range("D37").Value = '''42.42'
$ws.Range("E37").Value = '  +7.08%  '
$ws.Range("D38").Value = '0.0₃0869'
$ws.Range("E38").Value = '  +2.48%  '
$ws.Range("D39").Value = '''0.431'
$ws.Range("E39").Value = '  +0.10%  '
$ws.Range("D40").Value = '''0.152'
$ws.Range("E40").Value = '  +3.90%  '
$ws.Range("D41").Value = '''3.46'
$ws.Range("E41").Value = '  +8.01%  '
$ws.Range("E42").Value = '  +0.13%  '
$ws.Range("D43").Value = '''0.0501'
$ws.Range("E43").Value = '  +5.01%  '
$ws.Range("E44").Value = '  -0.14%  '
$ws.Range("E45").Value = '  +2.43%  '
$ws.Range("E46").Value = '  +13.68%  '
$ws.Range("E47").Value = '  -2.05%  '
$ws.Range("B48").Value = 'THORChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D48").Value = '''9.84'
$ws.Range("E48").Value = '  +16.72%  '
$ws.Range("B49").Value = 'ApeXProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D49").Value = '''3.42'
$ws.Range("E49").Value = '  +0.74%  '
$ws.Range("E50").Value = '  +5.19%  '
$ws.Range("E51").Value = '  +3.51%  '
